$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the cells whose values look like numbers or dates
# (e.g. "12345678", "03-01-1991", "08-13-2024") so Excel stores them as
# plain text instead of auto-converting to a number/date serial value,
# matching the text values in the target workbook.
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"

# Header row
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Custom Label"
$ws.Range("D1").Value = "Value"

# Data rows: Type, Custom Label, Value for rows 2-26
$data = @{
    2  = @("input", "NAME_OF_INDIVIDUAL", "Ashok Rohidas Modhave")
    3  = @("input", "COUNTRY_OF_RESIDENCE", "India")
    4  = @("input", "", "Aundh, Pune")
    5  = @("input", "", "Pune, Maharashtra 411011")
    6  = @("input", "", "India")
    7  = @("input", "", "")
    8  = @("input", "", "")
    9  = @("input", "", "")
    10 = @("input", "", "123-433-112")
    11 = @("input", "", "")
    12 = @("checkbox", "", "Unchecked")
    13 = @("input", "", "12345678")
    14 = @("input", "", "03-01-1991")
    15 = @("input", "", "India")
    16 = @("input", "", "")
    17 = @("input", "", "")
    18 = @("input", "", "")
    19 = @("input", "", "")
    20 = @("input", "", "")
    21 = @("checkbox", "", "Checked")
    22 = @("unknown", "", "")
    23 = @("input", "", "08-13-2024")
    24 = @("input", "", "Ashok Modhave")
    25 = @("unknown", "", "")
    26 = @("unknown", "", "")
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
}
